$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 750.1852
$ws.Range("J19").Value = 1023.2353
$ws.Range("L19").Value = 1023.2353
$ws.Range("N19").Value = -1373.2353
$ws.Range("H40").Value = 6667.2
$ws.Range("I40").Value = 6667.2
$ws.Range("K40").Value = 6667.2
$ws.Range("M40").Value = -6492.2
$ws.Range("H43").Value = 10580.059
$ws.Range("I43").Value = 17913.143
$ws.Range("K43").Value = 17913.143
$ws.Range("M43").Value = -17844.143
$ws.Range("H62").Value = 5744.5
$ws.Range("I62").Value = 5865.143
$ws.Range("J62").Value = 4900
$ws.Range("K62").Value = 5865.143
$ws.Range("L62").Value = 4900
$ws.Range("M62").Value = -5241.143
$ws.Range("N62").Value = -6148
$ws.Range("H65").Value = 5744.5
$ws.Range("I65").Value = 5865.143
$ws.Range("J65").Value = 4900
$ws.Range("K65").Value = 29325.715
$ws.Range("L65").Value = 24500
$ws.Range("M65").Value = -26205.715
$ws.Range("N65").Value = -30740
$ws.Range("H69").Value = 7999.3335
$ws.Range("J69").Value = 9749.5
$ws.Range("L69").Value = 29248.5
$ws.Range("N69").Value = -30996.5
$ws.Range("H72").Value = 7999.3335
$ws.Range("J72").Value = 9749.5
$ws.Range("L72").Value = 87745.5
$ws.Range("N72").Value = -96481.5
$ws.Range("H74").Value = 72979.734
$ws.Range("I74").Value = 81130.46000000001
$ws.Range("K74").Value = 81130.46000000001
$ws.Range("M74").Value = -80194.46000000001
$ws.Range("H77").Value = 72979.734
$ws.Range("I77").Value = 81130.46000000001
$ws.Range("K77").Value = 405652.3
$ws.Range("M77").Value = -400972.3
$ws.Range("H80").Value = 1971.2858
$ws.Range("I80").Value = 900
$ws.Range("J80").Value = 2399.8
$ws.Range("K80").Value = 2700
$ws.Range("L80").Value = 7199.400000000001
$ws.Range("M80").Value = -1702
$ws.Range("N80").Value = -9195.400000000001
$ws.Range("H83").Value = 1971.2858
$ws.Range("I83").Value = 900
$ws.Range("J83").Value = 2399.8
$ws.Range("K83").Value = 8100
$ws.Range("L83").Value = 21598.2
$ws.Range("M83").Value = -3108
$ws.Range("N83").Value = -31582.2
$ws.Range("H106").Value = 992
$ws.Range("I106").Value = 992
$ws.Range("K106").Value = 992
$ws.Range("M106").Value = -361
$ws.Range("H116").Value = 9448.5
$ws.Range("J116").Value = 9497.833000000001
$ws.Range("L116").Value = 9497.833000000001
$ws.Range("N116").Value = -16381.833
$ws.Range("H132").Value = 2074.8462
$ws.Range("I132").Value = 1957.84
$ws.Range("K132").Value = 5873.52
$ws.Range("M132").Value = -3343.52
$ws.Range("H137").Value = 2050.4
$ws.Range("I137").Value = 2475.375
$ws.Range("K137").Value = 7426.125
$ws.Range("M137").Value = -4876.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9245.088
$ws.Range("I32").Value = 4532.156
$ws.Range("K32").Value = 4532.156
$ws.Range("M32").Value = -4245.156
$ws.Range("H47").Value = 41000
$ws.Range("J47").Value = 41000
$ws.Range("L47").Value = 41000
$ws.Range("N47").Value = -42450
$ws.Range("H54").Value = 36833.332
$ws.Range("J54").Value = 36833.332
$ws.Range("L54").Value = 36833.332
$ws.Range("N54").Value = -38371.332
$ws.Range("H63").Value = 3710.6365
$ws.Range("J63").Value = 3581.7
$ws.Range("L63").Value = 3581.7
$ws.Range("N63").Value = -4953.7
$ws.Range("H66").Value = 3710.6365
$ws.Range("J66").Value = 3581.7
$ws.Range("L66").Value = 17908.5
$ws.Range("N66").Value = -24772.5
$ws.Range("H97").Value = 47136.863
$ws.Range("I97").Value = 1240.5555
$ws.Range("K97").Value = 1240.5555
$ws.Range("M97").Value = -744.5554999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2623.6
$ws.Range("I86").Value = 1714.1818
$ws.Range("K86").Value = 1714.1818
$ws.Range("M86").Value = -591.1818000000001
$ws.Range("H89").Value = 2623.6
$ws.Range("I89").Value = 1714.1818
$ws.Range("K89").Value = 8570.909
$ws.Range("M89").Value = -2954.909
$ws.Range("H134").Value = 2147.2188
$ws.Range("I134").Value = 2109.889
$ws.Range("J134").Value = 4499
$ws.Range("K134").Value = 6329.667
$ws.Range("L134").Value = 13497
$ws.Range("M134").Value = -3794.667
$ws.Range("N134").Value = -18567

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 69119.266
$ws.Range("I31").Value = 85085
$ws.Range("K31").Value = 85085
$ws.Range("M31").Value = -84790
$ws.Range("H34").Value = 69119.266
$ws.Range("I34").Value = 85085
$ws.Range("K34").Value = 85085
$ws.Range("M34").Value = -84883
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H60").Value = 16147.889
$ws.Range("J60").Value = 19995.8
$ws.Range("L60").Value = 19995.8
$ws.Range("N60").Value = -21017.8
$ws.Range("H62").Value = 5645.9375
$ws.Range("I62").Value = 3816.25
$ws.Range("J62").Value = 7475.625
$ws.Range("K62").Value = 3816.25
$ws.Range("L62").Value = 7475.625
$ws.Range("M62").Value = -3192.25
$ws.Range("N62").Value = -8723.625
$ws.Range("H65").Value = 5645.9375
$ws.Range("I65").Value = 3816.25
$ws.Range("J65").Value = 7475.625
$ws.Range("K65").Value = 19081.25
$ws.Range("L65").Value = 37378.125
$ws.Range("M65").Value = -15961.25
$ws.Range("N65").Value = -43618.125
$ws.Range("H68").Value = 100000
$ws.Range("J68").Value = 100000
$ws.Range("L68").Value = 100000
$ws.Range("N68").Value = -101498
$ws.Range("H71").Value = 100000
$ws.Range("J71").Value = 100000
$ws.Range("L71").Value = 300000
$ws.Range("N71").Value = -307488

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 46000
$ws.Range("J53").Value = 46000
$ws.Range("L53").Value = 46000
$ws.Range("N53").Value = -47262

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2649.8333
$ws.Range("I46").Value = 1299.6666
$ws.Range("K46").Value = 1299.6666
$ws.Range("M46").Value = -1111.6666
$ws.Range("H55").Value = 179.55556
$ws.Range("I55").Value = 154.75
$ws.Range("K55").Value = 154.75
$ws.Range("M55").Value = 18.25
$ws.Range("H57").Value = 36250
$ws.Range("J57").Value = 36250
$ws.Range("L57").Value = 36250
$ws.Range("N57").Value = -37382
$ws.Range("H68").Value = 627952.4
$ws.Range("I68").Value = 2802.7856
$ws.Range("K68").Value = 2802.7856
$ws.Range("M68").Value = -2053.7856
$ws.Range("H71").Value = 627952.4
$ws.Range("I71").Value = 2802.7856
$ws.Range("K71").Value = 14013.928
$ws.Range("M71").Value = -10269.928
$ws.Range("H100").Value = 18842
$ws.Range("I100").Value = 13649.167
$ws.Range("K100").Value = 13649.167
$ws.Range("M100").Value = -13108.167
$ws.Range("H140").Value = 147307.83
$ws.Range("I140").Value = 69949.336
$ws.Range("J140").Value = 224666.33
$ws.Range("K140").Value = 69949.336
$ws.Range("L140").Value = 224666.33
$ws.Range("M140").Value = -64769.336
$ws.Range("N140").Value = -235026.33

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2069.889
$ws.Range("I126").Value = 1988.5
$ws.Range("J126").Value = 2232.6667
$ws.Range("K126").Value = 5965.5
$ws.Range("L126").Value = 6698.000100000001
$ws.Range("M126").Value = -3495.5
$ws.Range("N126").Value = -11638.0001
